$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49 (shifts existing rows 49-59 down to 50-60),
# making room for a new weekly price entry at the top of this block.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly record.
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).Value = 44855
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = 100112031
$ws.Cells.Item(49, 7).Value = "Poroto verde"
$ws.Cells.Item(49, 8).Value = "Magnum"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 100
$ws.Cells.Item(49, 11).Value = 30000
$ws.Cells.Item(49, 12).Value = 32000
$ws.Cells.Item(49, 13).Value = 31000
$ws.Cells.Item(49, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(49, 15).Value = "Perú"
$ws.Cells.Item(49, 16).Value = 1240
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = "Hortaliza"
